# Fruta / hortaliza, semanal
# Insert a new weekly observation row above row 75 (shifting the existing
# rows 75-110 down to 76-111) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 75; everything below (old rows 75-110)
# shifts down by one (becoming rows 76-111).
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with this week's record.
$ws.Cells.Item(75, 1).Value = 8
$ws.Cells.Item(75, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(75, 3).Value = "Coquimbo"
$ws.Cells.Item(75, 4).Value = 44452
$ws.Cells.Item(75, 5).Value = 4
$ws.Cells.Item(75, 6).Value = 100112037
$ws.Cells.Item(75, 7).Value = "Cebollín"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 3100
$ws.Cells.Item(75, 11).Value = 900
$ws.Cells.Item(75, 12).Value = 1000
$ws.Cells.Item(75, 13).Value = 950
$ws.Cells.Item(75, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(75, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(75, 16).Value = 158
$ws.Cells.Item(75, 17).Value = 6
$ws.Cells.Item(75, 18).Value = "Hortaliza"
